# Scheduled-runner update: refresh cached market-board price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the per-job Leve sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ---- ALC -------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 18 - You Grow, Girl / Growth Formula Beta (prices now populated)
$ws.Range("H18").Value = 974.6667
$ws.Range("I18").Value = 769.6
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 769.6
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -485.6
$ws.Range("N18").Value = -2568

# Row 112 - Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 558444.9399999999
$ws.Range("J112").Value = 669834.9
$ws.Range("L112").Value = 2009504.7
$ws.Range("N112").Value = -2011720.7

# Row 138 - All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 154381.89
$ws.Range("I138").Value = 33210.484
$ws.Range("J138").Value = 245999.3
$ws.Range("K138").Value = 99631.45199999999
$ws.Range("L138").Value = 737997.8999999999
$ws.Range("M138").Value = -94491.45199999999
$ws.Range("N138").Value = -748277.8999999999

# ---- ARM -------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 24 - A Firm Hand / Iron Gauntlets
$ws.Range("H24").Value = 89000
$ws.Range("J24").Value = 89000
$ws.Range("L24").Value = 89000
$ws.Range("N24").Value = -89748

# Row 28 - 246 Kinds of Cheese / Iron Frypan
$ws.Range("H28").Value = 16006.2
$ws.Range("I28").Value = 5007.75
$ws.Range("J28").Value = 60000
$ws.Range("K28").Value = 5007.75
$ws.Range("L28").Value = 60000
$ws.Range("M28").Value = -4815.75
$ws.Range("N28").Value = -60384

# Row 31 - I Was a Teenage Wailer / Iron Alembic
$ws.Range("H31").Value = 978.25
$ws.Range("I31").Value = 978.25
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 978.25
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -684.25
$ws.Range("N31").ClearContents()

# Row 32 - Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3116.693
$ws.Range("I32").Value = 2983.465
$ws.Range("J32").Value = 8845.5
$ws.Range("K32").Value = 2983.465
$ws.Range("L32").Value = 8845.5
$ws.Range("M32").Value = -2696.465
$ws.Range("N32").Value = -9419.5

# Row 92 - Mail It In / High Steel Scale Mail of Fending
$ws.Range("H92").Value = 61094.5
$ws.Range("J92").Value = 61094.5
$ws.Range("L92").Value = 61094.5
$ws.Range("N92").Value = -66086.5

# Row 94 - Setting the Stage / High Steel Helm of Maiming
$ws.Range("H94").Value = 39015
$ws.Range("J94").Value = 39015
$ws.Range("L94").Value = 39015
$ws.Range("N94").Value = -40817

# Row 95 - Shielded Life / High Steel Scutum (prices cleared back to 0)
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 96 - The Gauntlet Is Cast / High Steel Gauntlets of Fending (cleared)
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# Row 98 - Greaving / Doman Iron Greaves of Maiming
$ws.Range("H98").Value = 58332.668
$ws.Range("J98").Value = 58332.668
$ws.Range("L98").Value = 58332.668
$ws.Range("N98").Value = -64322.668

# Row 99 - Home Cooking / Doman Iron Frypan
$ws.Range("H99").Value = 16006.2
$ws.Range("I99").Value = 5007.75
$ws.Range("J99").Value = 60000
$ws.Range("K99").Value = 5007.75
$ws.Range("L99").Value = 60000
$ws.Range("M99").Value = -2012.75
$ws.Range("N99").Value = -65990

# Row 100 - En Garde and on Guard / Doman Iron Gauntlets of Fending
$ws.Range("H100").Value = 89000
$ws.Range("J100").Value = 89000
$ws.Range("L100").Value = 89000
$ws.Range("N100").Value = -91164

# Row 101 - Art Imitates Life / Doman Steel Tabard of Fending
$ws.Range("H101").Value = 80999.664
$ws.Range("J101").Value = 80999.664
$ws.Range("L101").Value = 80999.664
$ws.Range("N101").Value = -87489.664

# Row 102 - Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 6669565
$ws.Range("I102").Value = 3025.28
$ws.Range("K102").Value = 3025.28
$ws.Range("M102").Value = -1403.28

# Row 103 - Sweeping the Legs / Doman Steel Greaves of Striking
$ws.Range("H103").Value = 100000
$ws.Range("J103").Value = 100000
$ws.Range("L103").Value = 100000
$ws.Range("N103").Value = -102344

# Row 104 - See Shields by the Sea Shore / Molybdenum Kite Shield
$ws.Range("H104").Value = 20112.5
$ws.Range("J104").Value = 20112.5
$ws.Range("L104").Value = 20112.5
$ws.Range("N104").Value = -27100.5

# Row 106 - Heads Will Roll / Molybdenum Headgear of Maiming
$ws.Range("H106").Value = 89443
$ws.Range("J106").Value = 89443
$ws.Range("L106").Value = 89443
$ws.Range("N106").Value = -91967

# ---- BSM -------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 86 - Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 22229192
$ws.Range("I86").Value = 22229192
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 22229192
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -22228069
$ws.Range("N86").ClearContents()

# Row 89 - Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 22229192
$ws.Range("I89").Value = 22229192
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 111145960
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -111140344
$ws.Range("N89").ClearContents()

# ---- CRP -------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 22 - Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 695.8182
$ws.Range("J22").Value = 687.5
$ws.Range("L22").Value = 687.5
$ws.Range("N22").Value = -1387.5

# Row 31 - Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2099.0984
$ws.Range("I31").Value = 1563.4546
$ws.Range("J31").Value = 3485.4707
$ws.Range("K31").Value = 1563.4546
$ws.Range("L31").Value = 3485.4707
$ws.Range("M31").Value = -1268.4546
$ws.Range("N31").Value = -4075.4707

# Row 34 - Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2099.0984
$ws.Range("I34").Value = 1563.4546
$ws.Range("J34").Value = 3485.4707
$ws.Range("K34").Value = 1563.4546
$ws.Range("L34").Value = 3485.4707
$ws.Range("M34").Value = -1361.4546
$ws.Range("N34").Value = -3889.4707

# Row 62 - Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 5555
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4931
$ws.Range("N62").ClearContents()

# Row 65 - The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 5555
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -24655
$ws.Range("N65").ClearContents()

# ---- CUL -------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 5 - What a Sap / Maple Syrup
$ws.Range("H5").Value = 7117.1763
$ws.Range("I5").Value = 449.5
$ws.Range("J5").Value = 8006.2
$ws.Range("K5").Value = 1348.5
$ws.Range("L5").Value = 24018.6
$ws.Range("M5").Value = -1236.5
$ws.Range("N5").Value = -24242.6

# Row 135 - Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 7117.1763
$ws.Range("I135").Value = 449.5
$ws.Range("J135").Value = 8006.2
$ws.Range("K135").Value = 4045.5
$ws.Range("L135").Value = 72055.8
$ws.Range("M135").Value = -1510.5
$ws.Range("N135").Value = -77125.8

# ---- LTW -------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 20 - Choke Hold / Hard Leather Choker
$ws.Range("H20").Value = 19979.334
$ws.Range("J20").Value = 19969
$ws.Range("L20").Value = 19969
$ws.Range("N20").Value = -20421

# Row 68 - You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 1795
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# Row 71 - They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 1795
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# ---- WVR -------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 25 - A Drag of a Doublet / Initiate's Doublet Vest
$ws.Range("H25").Value = 27500
$ws.Range("J25").Value = 27500
$ws.Range("L25").Value = 27500
$ws.Range("N25").Value = -28086

# Row 81 - Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 45307.04
$ws.Range("I81").Value = 115427.11
$ws.Range("J81").Value = 8184.647
$ws.Range("K81").Value = 230854.22
$ws.Range("L81").Value = 16369.294
$ws.Range("M81").Value = -229793.22
$ws.Range("N81").Value = -18491.294

# Row 84 - To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 45307.04
$ws.Range("I84").Value = 115427.11
$ws.Range("J84").Value = 8184.647
$ws.Range("K84").Value = 1154271.1
$ws.Range("L84").Value = 81846.47
$ws.Range("M84").Value = -1148967.1
$ws.Range("N84").Value = -92454.47

# Row 132 - Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2195.1667
$ws.Range("I132").Value = 2202.6316
$ws.Range("K132").Value = 6607.8948
$ws.Range("M132").Value = -4077.8948

Write-Output "Faerie Profits sheets updated."
